$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update Emilia's balance (row 2, column C) from 206183.69 to 201049.08
$ws.Cells.Item(2, 3).Value = 201049.08

# 2. Remove the rows for Rosane, Sofia, Cezar and Gabrielle (rows 4-7)
$ws.Range("A4:A7").EntireRow.Delete()

# 3. Insert a new row in their place with Patricia's new, smaller balance entry
$ws.Rows.Item(4).Insert()

# Force the account number into column A to be kept as text (preserve the
# leading zeros), then strip the temporary text format back off again so the
# cell ends up unformatted, like the rest of the data rows.
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "004421636"
$ws.Cells.Item(4, 1).ClearFormats()

$ws.Cells.Item(4, 2).Value = "Patricia"
$ws.Cells.Item(4, 3).Value = 10000

# 4. Remove the old duplicate Patricia row (previously row 12, now row 9
#    after the 4 deletions and 1 insertion above)
$ws.Range("A9:A9").EntireRow.Delete()
